$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of data (rows 9-12), matching the existing column layout:
# A: Date (serial), B: totalScore, C: posWordPercentage, D: negWordPercentage,
# E: posPhrasePercentage, F: negPhrasePercentage, G: ElapsedMs, H: wordCount,
# I: sentenceCount, J: posWordCount, K: negWordCount, L: positivePhraseCount,
# M: negativePhraseCount, N: Method ("Bag")

$rows = @(
    @{ r=9;  A=42613.761018518519; B=16; C=54; D=43; E=54; F=28; G=24747; H=27948; I=3135; J=419; K=334; L=28; M=11 },
    @{ r=10; A=42613.890694444446; B=14; C=51; D=42; E=51; F=31; G=12667; H=25961; I=2889; J=393; K=326; L=31; M=14 },
    @{ r=11; A=42614.88753472222;  B=34; C=60; D=37; E=60; F=13; G=18171; H=18158; I=2005; J=290; K=177; L=38; M=6  },
    @{ r=12; A=42615.886805555558; B=34; C=59; D=38; E=59; F=9;  G=11517; H=23194; I=2606; J=374; K=240; L=40; M=4  }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = "Bag"
}
